$d = $word.ActiveDocument

function Get-ParagraphAt($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -le $pos -and $p.Range.End -gt $pos) {
            return $p
        }
    }
    return $null
}

# --- Locate the "Architecture modele-vue-controleur" heading paragraph ---
$headingText = "Architecture mod" + [char]0x00E8 + "le-vue-contr" + [char]0x00F4 + "leur"
$headingRange = $d.Content.Duplicate
$headingRange.Find.Execute($headingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingPara = Get-ParagraphAt $d $headingRange.Start

# --- Remove the duplicate empty "page break only" paragraph directly in
#     front of that heading (it held a now-redundant lastRenderedPageBreak
#     run that Word had written into the previous page-break paragraph). ---
if ($headingPara -ne $null) {
    $prevPara = $headingPara.Previous()
    if ($prevPara -ne $null) {
        $prevText = $prevPara.Range.Text
        if ($prevText -eq ([string][char]12 + [string][char]13)) {
            $prevPara.Range.Delete()
        }
    }
}

# --- Re-resolve the heading after the deletion and move the document's
#     hidden "_GoBack" bookmark (Word always keeps exactly one, tracking the
#     last edit location) to the very start of the heading. Adding a
#     bookmark with an existing name relocates it instead of duplicating it. ---
$headingRange2 = $d.Content.Duplicate
$headingRange2.Find.Execute($headingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackRange = $headingRange2.Duplicate
$goBackRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $goBackRange)
